$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 12.01751833333333
$ws.Range("H2").Value = 36.052555
$ws.Range("I2").Value = 0.5714697501126412
$ws.Range("J2").Value = 0.6009086862203552
$ws.Range("M2").Value = 29.785352
$ws.Range("N2").Value = 89.356056
$ws.Range("O2").Value = 0.7923195065866085
$ws.Range("P2").Value = 0.7947519366640845
$ws.Range("Q2").Value = 357.9460137247867
$ws.Range("R2").Value = 3221.51412352308
$ws.Range("S2").Value = 0.4527866304384204
$ws.Range("T2").Value = 0.4775733421318979

# Row 3
$ws.Range("G3").Value = 12.01751833333333
$ws.Range("H3").Value = 36.052555
$ws.Range("I3").Value = 0.5714697501126412
$ws.Range("J3").Value = 0.6009086862203552
$ws.Range("O3").Value = 0.1171985110386058
$ws.Range("P3").Value = 0.1175583118271966
$ws.Range("Q3").Value = 52.946746220445
$ws.Range("R3").Value = 476.5207159840049
$ws.Range("S3").Value = 0.0669754038168057
$ws.Range("T3").Value = 0.07064181071436354

# Row 4
$ws.Range("G4").Value = 12.01751833333333
$ws.Range("H4").Value = 36.052555
$ws.Range("I4").Value = 0.5714697501126412
$ws.Range("J4").Value = 0.6009086862203552
$ws.Range("M4").Value = 1.362560333333333
$ws.Range("N4").Value = 4.087681
$ws.Range("O4").Value = 0.03624543805965938
$ws.Range("P4").Value = 0.03635671197501131
$ws.Range("Q4").Value = 16.37459378610611
$ws.Range("R4").Value = 147.371344074955
$ws.Range("S4").Value = 0.02071317143067676
$ws.Range("T4").Value = 0.0218470640281959

# Row 5
$ws.Range("G5").Value = 12.01751833333333
$ws.Range("H5").Value = 36.052555
$ws.Range("I5").Value = 0.5714697501126412
$ws.Range("J5").Value = 0.6009086862203552
$ws.Range("M5").Value = 0.3451695
$ws.Range("N5").Value = 0.690339
$ws.Range("O5").Value = 0.00918184643004207
$ws.Range("P5").Value = 0.006140023203404898
$ws.Range("Q5").Value = 4.1480807943575
$ws.Range("R5").Value = 24.888484766145
$ws.Range("S5").Value = 0.005247147484948788
$ws.Range("T5").Value = 0.003689593276520534

# Row 6
$ws.Range("G6").Value = 12.01751833333333
$ws.Range("H6").Value = 36.052555
$ws.Range("I6").Value = 0.5714697501126412
$ws.Range("J6").Value = 0.6009086862203552
$ws.Range("M6").Value = 1.693723333333333
$ws.Range("N6").Value = 5.08117
$ws.Range("O6").Value = 0.04505469788508434
$ws.Range("P6").Value = 0.04519301633030275
$ws.Range("Q6").Value = 20.35435120992778
$ws.Range("R6").Value = 183.18916088935
$ws.Range("S6").Value = 0.02574739694178969
$ws.Range("T6").Value = 0.02715687606937728

# Row 7
$ws.Range("I7").Value = 0.03126852934637925
$ws.Range("J7").Value = 0.03287930968502149
$ws.Range("M7").Value = 29.785352
$ws.Range("N7").Value = 89.356056
$ws.Range("O7").Value = 0.7923195065866085
$ws.Range("P7").Value = 0.7947519366640845
$ws.Range("Q7").Value = 19.58536813605067
$ws.Range("R7").Value = 176.268313224456
$ws.Range("S7").Value = 0.02477466574341209
$ws.Range("T7").Value = 0.02613089504834902

# Row 8
$ws.Range("I8").Value = 0.03126852934637925
$ws.Range("J8").Value = 0.03287930968502149
$ws.Range("O8").Value = 0.1171985110386058
$ws.Range("P8").Value = 0.1175583118271966
$ws.Range("S8").Value = 0.003664625081762599
$ws.Range("T8").Value = 0.003865236140614721

# Row 9
$ws.Range("I9").Value = 0.03126852934637925
$ws.Range("J9").Value = 0.03287930968502149
$ws.Range("M9").Value = 1.362560333333333
$ws.Range("N9").Value = 4.087681
$ws.Range("O9").Value = 0.03624543805965938
$ws.Range("P9").Value = 0.03635671197501131
$ws.Range("Q9").Value = 0.8959520013701111
$ws.Range("R9").Value = 8.063568012331
$ws.Range("S9").Value = 0.001133341543640831
$ws.Range("T9").Value = 0.001195383592155526

# Row 10
$ws.Range("I10").Value = 0.03126852934637925
$ws.Range("J10").Value = 0.03287930968502149
$ws.Range("M10").Value = 0.3451695
$ws.Range("N10").Value = 0.690339
$ws.Range("O10").Value = 0.00918184643004207
$ws.Range("P10").Value = 0.006140023203404898
$ws.Range("Q10").Value = 0.2269663197815
$ws.Range("R10").Value = 1.361797918689
$ws.Range("S10").Value = 0.000287102834551718
$ws.Range("T10").Value = 0.0002018797243779674

# Row 11
$ws.Range("I11").Value = 0.03126852934637925
$ws.Range("J11").Value = 0.03287930968502149
$ws.Range("M11").Value = 1.693723333333333
$ws.Range("N11").Value = 5.08117
$ws.Range("O11").Value = 0.04505469788508434
$ws.Range("P11").Value = 0.04519301633030275
$ws.Range("Q11").Value = 1.113708342407778
$ws.Range("R11").Value = 10.02337508167
$ws.Range("S11").Value = 0.001408794143012011
$ws.Range("T11").Value = 0.001485915179524258

# Row 12
$ws.Range("G12").Value = 3.025265666666666
$ws.Range("H12").Value = 9.075797
$ws.Range("I12").Value = 0.1438606346668928
$ws.Range("J12").Value = 0.1512715326742485
$ws.Range("M12").Value = 29.785352
$ws.Range("N12").Value = 89.356056
$ws.Range("O12").Value = 0.7923195065866085
$ws.Range("P12").Value = 0.7947519366640845
$ws.Range("Q12").Value = 90.10860277518132
$ws.Range("R12").Value = 810.9774249766319
$ws.Range("S12").Value = 0.1139835870765088
$ws.Range("T12").Value = 0.1202233435550033

# Row 13
$ws.Range("G13").Value = 3.025265666666666
$ws.Range("H13").Value = 9.075797
$ws.Range("I13").Value = 0.1438606346668928
$ws.Range("J13").Value = 0.1512715326742485
$ws.Range("O13").Value = 0.1171985110386058
$ws.Range("P13").Value = 0.1175583118271966
$ws.Range("Q13").Value = 13.328706398403
$ws.Range("R13").Value = 119.958357585627
$ws.Range("S13").Value = 0.01686025218002868
$ws.Range("T13").Value = 0.01778322600869726

# Row 14
$ws.Range("G14").Value = 3.025265666666666
$ws.Range("H14").Value = 9.075797
$ws.Range("I14").Value = 0.1438606346668928
$ws.Range("J14").Value = 0.1512715326742485
$ws.Range("M14").Value = 1.362560333333333
$ws.Range("N14").Value = 4.087681
$ws.Range("O14").Value = 0.03624543805965938
$ws.Range("P14").Value = 0.03635671197501131
$ws.Range("Q14").Value = 4.122106995195222
$ws.Range("R14").Value = 37.098962956757
$ws.Range("S14").Value = 0.005214291723042149
$ws.Range("T14").Value = 0.005499735543456166

# Row 15
$ws.Range("G15").Value = 3.025265666666666
$ws.Range("H15").Value = 9.075797
$ws.Range("I15").Value = 0.1438606346668928
$ws.Range("J15").Value = 0.1512715326742485
$ws.Range("M15").Value = 0.3451695
$ws.Range("N15").Value = 0.690339
$ws.Range("O15").Value = 0.00918184643004207
$ws.Range("P15").Value = 0.006140023203404898
$ws.Range("Q15").Value = 1.0442294375305
$ws.Range("R15").Value = 6.265376625183
$ws.Range("S15").Value = 0.001320906254839796
$ws.Range("T15").Value = 0.0009288107206345081

# Row 16
$ws.Range("G16").Value = 3.025265666666666
$ws.Range("H16").Value = 9.075797
$ws.Range("I16").Value = 0.1438606346668928
$ws.Range("J16").Value = 0.1512715326742485
$ws.Range("M16").Value = 1.693723333333333
$ws.Range("N16").Value = 5.08117
$ws.Range("O16").Value = 0.04505469788508434
$ws.Range("P16").Value = 0.04519301633030275
$ws.Range("Q16").Value = 5.123963049165556
$ws.Range("R16").Value = 46.11566744249
$ws.Range("S16").Value = 0.006481597432473346
$ws.Range("T16").Value = 0.006836416846457238

# Row 17
$ws.Range("G17").Value = 3.090697
$ws.Range("H17").Value = 6.181394
$ws.Range("I17").Value = 0.146972094676554
$ws.Range("J17").Value = 0.1030288518400537
$ws.Range("M17").Value = 29.785352
$ws.Range("N17").Value = 89.356056
$ws.Range("O17").Value = 0.7923195065866085
$ws.Range("P17").Value = 0.7947519366640845
$ws.Range("Q17").Value = 92.057498070344
$ws.Range("R17").Value = 552.3449884220639
$ws.Range("S17").Value = 0.1164488575361276
$ws.Range("T17").Value = 0.0818823795321597

# Row 18
$ws.Range("G18").Value = 3.090697
$ws.Range("H18").Value = 6.181394
$ws.Range("I18").Value = 0.146972094676554
$ws.Range("J18").Value = 0.1030288518400537
$ws.Range("O18").Value = 0.1171985110386058
$ws.Range("P18").Value = 0.1175583118271966
$ws.Range("Q18").Value = 13.616983570509
$ws.Range("R18").Value = 81.701901423054
$ws.Range("S18").Value = 0.01722491066031714
$ws.Range("T18").Value = 0.01211189789181107

# Row 19
$ws.Range("G19").Value = 3.090697
$ws.Range("H19").Value = 6.181394
$ws.Range("I19").Value = 0.146972094676554
$ws.Range("J19").Value = 0.1030288518400537
$ws.Range("M19").Value = 1.362560333333333
$ws.Range("N19").Value = 4.087681
$ws.Range("O19").Value = 0.03624543805965938
$ws.Range("P19").Value = 0.03635671197501131
$ws.Range("Q19").Value = 4.211261134552333
$ws.Range("R19").Value = 25.267566807314
$ws.Range("S19").Value = 0.005327067954097433
$ws.Range("T19").Value = 0.003745790291464947

# Row 20
$ws.Range("G20").Value = 3.090697
$ws.Range("H20").Value = 6.181394
$ws.Range("I20").Value = 0.146972094676554
$ws.Range("J20").Value = 0.1030288518400537
$ws.Range("M20").Value = 0.3451695
$ws.Range("N20").Value = 0.690339
$ws.Range("O20").Value = 0.00918184643004207
$ws.Range("P20").Value = 0.006140023203404898
$ws.Range("Q20").Value = 1.0668143381415
$ws.Range("R20").Value = 4.267257352566
$ws.Range("S20").Value = 0.001349475202821722
$ws.Range("T20").Value = 0.0006325995409180951

# Row 21
$ws.Range("G21").Value = 3.090697
$ws.Range("H21").Value = 6.181394
$ws.Range("I21").Value = 0.146972094676554
$ws.Range("J21").Value = 0.1030288518400537
$ws.Range("M21").Value = 1.693723333333333
$ws.Range("N21").Value = 5.08117
$ws.Range("O21").Value = 0.04505469788508434
$ws.Range("P21").Value = 0.04519301633030275
$ws.Range("Q21").Value = 5.234785625163334
$ws.Range("R21").Value = 31.40871375098
$ws.Range("S21").Value = 0.006621783323190153
$ws.Range("T21").Value = 0.004656184583699889

# Row 22
$ws.Range("G22").Value = 2.238110333333333
$ws.Range("H22").Value = 6.714331
$ws.Range("I22").Value = 0.1064289911975326
$ws.Range("J22").Value = 0.1119116195803211
$ws.Range("M22").Value = 29.785352
$ws.Range("N22").Value = 89.356056
$ws.Range("O22").Value = 0.7923195065866085
$ws.Range("P22").Value = 0.7947519366640845
$ws.Range("Q22").Value = 66.66290409317067
$ws.Range("R22").Value = 599.9661368385359
$ws.Range("S22").Value = 0.08432576579213955
$ws.Range("T22").Value = 0.08894197639667449

# Row 23
$ws.Range("G23").Value = 2.238110333333333
$ws.Range("H23").Value = 6.714331
$ws.Range("I23").Value = 0.1064289911975326
$ws.Range("J23").Value = 0.1119116195803211
$ws.Range("O23").Value = 0.1171985110386058
$ws.Range("P23").Value = 0.1175583118271966
$ws.Range("Q23").Value = 9.860659792268999
$ws.Range("R23").Value = 88.745938130421
$ws.Range("S23").Value = 0.01247331929969171
$ws.Range("T23").Value = 0.01315614107170999

# Row 24
$ws.Range("G24").Value = 2.238110333333333
$ws.Range("H24").Value = 6.714331
$ws.Range("I24").Value = 0.1064289911975326
$ws.Range("J24").Value = 0.1119116195803211
$ws.Range("M24").Value = 1.362560333333333
$ws.Range("N24").Value = 4.087681
$ws.Range("O24").Value = 0.03624543805965938
$ws.Range("P24").Value = 0.03635671197501131
$ws.Range("Q24").Value = 3.049560361823444
$ws.Range("R24").Value = 27.446043256411
$ws.Range("S24").Value = 0.003857565408202202
$ws.Range("T24").Value = 0.004068738519738772

# Row 25
$ws.Range("G25").Value = 2.238110333333333
$ws.Range("H25").Value = 6.714331
$ws.Range("I25").Value = 0.1064289911975326
$ws.Range("J25").Value = 0.1119116195803211
$ws.Range("M25").Value = 0.3451695
$ws.Range("N25").Value = 0.690339
$ws.Range("O25").Value = 0.00918184643004207
$ws.Range("P25").Value = 0.006140023203404898
$ws.Range("Q25").Value = 0.7725274247015
$ws.Range("R25").Value = 4.635164548209
$ws.Range("S25").Value = 0.0009772146528800438
$ws.Range("T25").Value = 0.0006871399409537936

# Row 26
$ws.Range("G26").Value = 2.238110333333333
$ws.Range("H26").Value = 6.714331
$ws.Range("I26").Value = 0.1064289911975326
$ws.Range("J26").Value = 0.1119116195803211
$ws.Range("M26").Value = 1.693723333333333
$ws.Range("N26").Value = 5.08117
$ws.Range("O26").Value = 0.04505469788508434
$ws.Range("P26").Value = 0.04519301633030275
$ws.Range("Q26").Value = 3.790739694141111
$ws.Range("R26").Value = 34.11665724727
$ws.Range("S26").Value = 0.004795126044619133
$ws.Range("T26").Value = 0.005057623651244081
